$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New run data rows (rows 2-11)
$data = @(
    @(1, 69771, 6610, 34.877906799316413, 18977, 100),
    @(2, 69342, 8481, 45.788062810897827, 18801, 100),
    @(3, 190565, 290, 15.37570357322693, 1999, 100),
    @(4, 188241, 255, 13.51781487464905, 1618, 100),
    @(5, 365438, 486, 53.930134534835823, 892, 100),
    @(6, 74160, 73, 3.995342493057251, 2116, 100),
    @(7, 73087, 79, 5.7033231258392334, 1717, 100),
    @(8, 187428, 281, 25.20179200172424, 605, 100),
    @(9, 192348, 211, 28.320268630981449, 1032, 100),
    @(10, 367316, 396, 74.687691450119019, 514, 100)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}

# Update selection to match the final state
$ws.Range("C11").Select()
